$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 20230921
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 8

# Row 7 (row 6 intentionally left blank)
$ws.Range("A7").Value = 20230929
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 7

# Row 8 - shuffled order label
$ws.Range("A8").Value = "1,2,3,1,2,3"

# Row 9
$ws.Range("A9").Value = 20231005
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 6

# Row 10 - shuffled order label
$ws.Range("A10").Value = "2,3,1,2,3,1"

# Row 11
$ws.Range("A11").Value = 20231030
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 6

# Row 12 - shuffled order label (same text as row 10, reuses shared string)
$ws.Range("A12").Value = "2,3,1,2,3,1"

# Match the saved selection/active cell state from the diff
[void]$ws.Range("A12").Select()

Write-Output "edit applied"
